# Compute X Interpolants once.
# Adds a "X Delta once" row to the existing "Tex, Flags<0>" sheet with a
# percentage-delta formula, and adds a new "Flat, Flags<0>" sheet that
# mirrors the same Baseline layout (values only, delta left for later).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet 1 ("Tex, Flags<0>"): new row 3 with the X-delta formula -------
$ws1.Range("A3").Value = "X Delta once"
$ws1.Range("B3").Value = 71957
$ws1.Range("C3").Value = 73067
$ws1.Range("D3").Value = 76

$ws1.Range("E1").Value = "%"
$ws1.Range("E1").NumberFormat = "0.0%"
$ws1.Range("E3").NumberFormat = "0.0%"
$ws1.Range("E3").Formula = "=(D3/D2)-1"

# Column A needs to fit "X Delta once" -> 12 characters wide.
$ws1.Columns("A").ColumnWidth = 12 - (5 / 6)

[void]$ws1.Range("A3").Select()

# --- New sheet "Flat, Flags<0>" -------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Flat, Flags<0>"

$ws2.Range("B1").Value = "GCC Win"
$ws2.Range("C1").Value = "MSVC Win"
$ws2.Range("D1").Value = "GBA"
$ws2.Range("E1").Value = "%"

$ws2.Range("A2").Value = "Baseline"
$ws2.Range("B2").Value = 155448
$ws2.Range("C2").Value = 165920
$ws2.Range("D2").Value = 396

$ws2.Range("E1:E3").NumberFormat = "0.0%"

[void]$ws2.Range("A3").Select()

# Leave the original sheet as the active / selected tab.
$ws1.Activate()
